$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("candybars")
$dst = $wb.Worksheets.Item("candybars_mini")

$srcRange = $src.Range("A1:K26")
$dstRange = $dst.Range("A1:K26")
$dstRange.Value2 = $srcRange.Value2

$dst.Range("F1").Select()
